$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "39.955.53"
$ws.Range("D3").Value = "2.220.10"
$ws.Range("E3").Value = "  -0.75%  "
$ws.Range("E4").Value = "  -0.07%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "292.47"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -0.38%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "86.93"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -2.06%  "
$ws.Range("E7").Value = "  -0.91%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("E9").Value = "  -2.22%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "30.47"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -2.39%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "50.36"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +5.52%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.0780"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -1.44%  "
$ws.Range("E13").Value = "  +3.14%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "6.44"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -0.50%  "
$ws.Range("D15").Value = "2.562.15"
$ws.Range("E15").Value = "  -0.71%  "
$ws.Range("E16").Value = "  -3.08%  "
$ws.Range("D17").Value = "2.256.99"
$ws.Range("E17").Value = "  -0.96%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "0.733"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -0.85%  "
$ws.Range("D19").Value = "39.871.24"
$ws.Range("E19").Value = "  -0.65%  "
$ws.Range("D20").Value = "0.0₃0886"
$ws.Range("E20").Value = "  -0.49%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "11.13"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -4.73%  "
$ws.Range("E22").Value = "  -2.22%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "65.55"
$c.Style = "Normal"
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "237.03"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +0.05%  "
$ws.Range("E25").Value = "  +0.07%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "2.46"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -1.03%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "1.83"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -1.62%  "
$ws.Range("E28").Value = "  +7.50%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "23.26"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +1.01%  "
$ws.Range("E30").Value = "  -1.30%  "
$ws.Range("E31").Value = "  +3.16%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "31.81"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -3.95%  "
$ws.Range("E33").Value = "  +0.03%  "
$ws.Range("E34").Value = "  -0.61%  "
$ws.Range("E35").Value = "  +3.61%  "
$ws.Range("E36").Value = "  -1.39%  "
$ws.Range("E37").Value = "  -1.40%  "
$ws.Range("E38").Value = "  -0.72%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.0987"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -1.54%  "
$ws.Range("E40").Value = "  -0.15%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "15.20"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -6.82%  "
$ws.Range("D42").Value = "2.087.44"
$ws.Range("E42").Value = "  -1.42%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "3.71"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -3.88%  "
$ws.Range("E44").Value = "  -0.50%  "
$ws.Range("E45").Value = "  -2.71%  "
$ws.Range("E46").Value = "  -3.21%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "1.99"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -9.00%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "2.70"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +0.39%  "
$ws.Range("D49").Value = "2.433.93"
$ws.Range("E49").Value = "  -0.66%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "1.46"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -1.47%  "
$ws.Range("E51").Value = "  +1.87%  "
